$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column B (StatQuery) between the existing "query" column (A)
# and the "dbExcel"/"WebExcel" columns (old B/C, which become C/D). This
# shifts the old B/C columns right while preserving their stored widths and
# content, matching the sheet1.xml diff (dimension A1:C2 -> A1:D2).
$ws.Columns("B").Insert()

# --- Header row (row 1) ---
$ws.Range("B1").Value = "StatQuery"

# --- New query row (row 2). Column B inherited the wrap-text format from
# column A (s="1") automatically via the column Insert() above, matching
# the existing query cell A2's style. ---
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Other']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# --- Column B should have the same width as column A (75.81640625 chars) ---
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth()

# --- Reset the view: scroll back to the top-left and select A2 (was scrolled
# to A2 with B3:B14 selected) ---
[void]$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A2").Select()
